$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 674.06665
$ws.Range("I11").Value = 674.06665
$ws.Range("K11").Value = 674.06665
$ws.Range("M11").Value = -534.06665
$ws.Range("H18").Value = 1995
$ws.Range("I18").Value = 1995
$ws.Range("K18").Value = 1995
$ws.Range("M18").Value = -1711
$ws.Range("H53").Value = 1183.3334
$ws.Range("I53").Value = 500
$ws.Range("J53").Value = 1525
$ws.Range("K53").Value = 500
$ws.Range("L53").Value = 1525
$ws.Range("M53").Value = 137
$ws.Range("N53").Value = -2799
$ws.Range("H103").Value = 1550.0769
$ws.Range("I103").Value = 1818.3334
$ws.Range("J103").Value = 1320.1428
$ws.Range("K103").Value = 5455.0002
$ws.Range("L103").Value = 3960.4284
$ws.Range("M103").Value = -4869.0002
$ws.Range("N103").Value = -5132.428400000001
$ws.Range("H109").Value = 57665.223
$ws.Range("J109").Value = 57665.223
$ws.Range("L109").Value = 57665.223
$ws.Range("N109").Value = -60439.223
$ws.Range("H110").Value = 67850.664
$ws.Range("J110").Value = 67850.664
$ws.Range("L110").Value = 67850.664
$ws.Range("N110").Value = -76030.664
$ws.Range("H111").Value = 899.8570999999999
$ws.Range("J111").Value = 1099.6666
$ws.Range("L111").Value = 3298.9998
$ws.Range("N111").Value = -9432.9998
$ws.Range("H113").Value = 4099.5
$ws.Range("J113").Value = 4459.4
$ws.Range("L113").Value = 4459.4
$ws.Range("N113").Value = -10967.4
$ws.Range("H116").Value = 3708445.2
$ws.Range("J116").Value = 6670321.5
$ws.Range("L116").Value = 6670321.5
$ws.Range("N116").Value = -6677205.5
$ws.Range("H137").Value = 520536.44
$ws.Range("I137").Value = 1966.3529
$ws.Range("J137").Value = 1321962.9
$ws.Range("K137").Value = 5899.0587
$ws.Range("L137").Value = 3965888.7
$ws.Range("M137").Value = -3349.0587
$ws.Range("N137").Value = -3970988.7
$ws.Range("H138").Value = 1239.2941
$ws.Range("J138").Value = 2500
$ws.Range("L138").Value = 7500
$ws.Range("N138").Value = -17780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6832.1772
$ws.Range("I32").Value = 2395.9268
$ws.Range("J32").Value = 15493.429
$ws.Range("K32").Value = 2395.9268
$ws.Range("L32").Value = 15493.429
$ws.Range("M32").Value = -2108.9268
$ws.Range("N32").Value = -16067.429
$ws.Range("H52").Value = 56558.4
$ws.Range("J52").Value = 56558.4
$ws.Range("L52").Value = 56558.4
$ws.Range("N52").Value = -57194.4
$ws.Range("H74").Value = 50553.332
$ws.Range("I74").Value = 85235.664
$ws.Range("K74").Value = 85235.664
$ws.Range("M74").Value = -84361.664
$ws.Range("H77").Value = 50553.332
$ws.Range("I77").Value = 85235.664
$ws.Range("K77").Value = 426178.32
$ws.Range("M77").Value = -421810.32
$ws.Range("H97").Value = 1575.875
$ws.Range("I97").Value = 1322
$ws.Range("J97").Value = 1999
$ws.Range("K97").Value = 1322
$ws.Range("L97").Value = 1999
$ws.Range("M97").Value = -826
$ws.Range("N97").Value = -2991
$ws.Range("H104").Value = 18899
$ws.Range("J104").Value = 20368.5
$ws.Range("L104").Value = 20368.5
$ws.Range("N104").Value = -27356.5
$ws.Range("H110").Value = 704.2
$ws.Range("I110").Value = 526.7646999999999
$ws.Range("K110").Value = 526.7646999999999
$ws.Range("M110").Value = 1518.2353
$ws.Range("H121").Value = 81487.125
$ws.Range("J121").Value = 81487.125
$ws.Range("L121").Value = 81487.125
$ws.Range("N121").Value = -84981.125
$ws.Range("H141").Value = 149961.75
$ws.Range("J141").Value = 149961.75
$ws.Range("L141").Value = 149961.75
$ws.Range("N141").Value = -160321.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 31100
$ws.Range("J2").Value = 31100
$ws.Range("L2").Value = 31100
$ws.Range("N2").Value = -31326
$ws.Range("H36").Value = 5972
$ws.Range("I36").Value = 1296
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 1296
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = -762
$ws.Range("N36").Value = -21068
$ws.Range("H94").Value = 2004.2693
$ws.Range("J94").Value = 3765.6667
$ws.Range("L94").Value = 3765.6667
$ws.Range("N94").Value = -4667.6667
$ws.Range("H109").Value = 78282.71000000001
$ws.Range("J109").Value = 78282.71000000001
$ws.Range("L109").Value = 78282.71000000001
$ws.Range("N109").Value = -81056.71000000001
$ws.Range("H112").Value = 74990
$ws.Range("J112").Value = 74990
$ws.Range("L112").Value = 74990
$ws.Range("N112").Value = -77944
$ws.Range("H119").Value = 98996
$ws.Range("J119").Value = 98996
$ws.Range("L119").Value = 98996
$ws.Range("N119").Value = -108672

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 94996
$ws.Range("J9").Value = 94996
$ws.Range("L9").Value = 94996
$ws.Range("N9").Value = -95332
$ws.Range("H22").Value = 399.6
$ws.Range("I22").Value = 249.5
$ws.Range("K22").Value = 249.5
$ws.Range("M22").Value = 100.5
$ws.Range("H31").Value = 2783.0967
$ws.Range("I31").Value = 1773.826
$ws.Range("J31").Value = 5684.75
$ws.Range("K31").Value = 1773.826
$ws.Range("L31").Value = 5684.75
$ws.Range("M31").Value = -1478.826
$ws.Range("N31").Value = -6274.75
$ws.Range("H34").Value = 2783.0967
$ws.Range("I34").Value = 1773.826
$ws.Range("J34").Value = 5684.75
$ws.Range("K34").Value = 1773.826
$ws.Range("L34").Value = 5684.75
$ws.Range("M34").Value = -1571.826
$ws.Range("N34").Value = -6088.75
$ws.Range("H108").Value = 44940.285
$ws.Range("J108").Value = 44940.285
$ws.Range("L108").Value = 44940.285
$ws.Range("N108").Value = -52620.285
$ws.Range("H114").Value = 39746
$ws.Range("J114").Value = 39746
$ws.Range("L114").Value = 39746
$ws.Range("N114").Value = -48424
$ws.Range("H116").Value = 65898.5
$ws.Range("J116").Value = 65898.5
$ws.Range("L116").Value = 65898.5
$ws.Range("N116").Value = -75076.5
$ws.Range("H119").Value = 99998.2
$ws.Range("J119").Value = 99998.2
$ws.Range("L119").Value = 99998.2
$ws.Range("N119").Value = -109674.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 6860.75
$ws.Range("J125").Value = 15000
$ws.Range("L125").Value = 45000
$ws.Range("N125").Value = -54840
$ws.Range("H128").Value = 357824.62
$ws.Range("I128").Value = 357824.62
$ws.Range("K128").Value = 1073473.86
$ws.Range("M128").Value = -1068493.86
$ws.Range("H131").Value = 1448.1818
$ws.Range("J131").Value = 1522.6666
$ws.Range("L131").Value = 4567.9998
$ws.Range("N131").Value = -14647.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 442.70834
$ws.Range("I107").Value = 370.3125
$ws.Range("J107").Value = 587.5
$ws.Range("K107").Value = 370.3125
$ws.Range("L107").Value = 587.5
$ws.Range("M107").Value = 1549.6875
$ws.Range("N107").Value = -4427.5
$ws.Range("H108").Value = 43994
$ws.Range("J108").Value = 43994
$ws.Range("L108").Value = 43994
$ws.Range("N108").Value = -51674
$ws.Range("H114").Value = 76530
$ws.Range("J114").Value = 76530
$ws.Range("L114").Value = 76530
$ws.Range("N114").Value = -85208
$ws.Range("H132").Value = 3964.7896
$ws.Range("I132").Value = 2758.0908
$ws.Range("K132").Value = 8274.2724
$ws.Range("M132").Value = -5744.2724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7483.25
$ws.Range("I22").Value = 1411
$ws.Range("J22").Value = 25700
$ws.Range("K22").Value = 1411
$ws.Range("L22").Value = 25700
$ws.Range("M22").Value = -1116
$ws.Range("N22").Value = -26290
$ws.Range("H27").Value = 7483.25
$ws.Range("I27").Value = 1411
$ws.Range("J27").Value = 25700
$ws.Range("K27").Value = 1411
$ws.Range("L27").Value = 25700
$ws.Range("M27").Value = -1304
$ws.Range("N27").Value = -25914
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H123").Value = 72169.875
$ws.Range("J123").Value = 75281.28999999999
$ws.Range("L123").Value = 75281.28999999999
$ws.Range("N123").Value = -85081.28999999999
$ws.Range("H136").Value = 3305.2173
$ws.Range("I136").Value = 3305.2173
$ws.Range("K136").Value = 9915.651899999999
$ws.Range("M136").Value = -7365.651899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1825.3478
$ws.Range("I132").Value = 1528.4375
$ws.Range("K132").Value = 4585.3125
$ws.Range("M132").Value = -2055.3125

Write-Output "done"